# Apply the LinuxForHealth re-brand + regen-refresh edit to the
# "StructureDefinition-hierarchy-version-date" FHIR IG spreadsheet.
#
# Changes (per the commit's xml diff):
#   Metadata sheet:
#     B2 (URL)       ibm.com -> linuxforhealth.org
#     B3 (Version)   7.0.0   -> 8.0.0
#     B8 (Date)      2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
#     B9 (Publisher) Alvearie Team -> LinuxForHealth Team
#   Elements sheet:
#     AI2 (Constraint(s) for the root "Extension" row) cleared to blank -
#     the ele-1/ext-1 constraint text now only lives on the
#     Extension.extension row (AI4), not on the Extension row itself.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/hierarchy-version-date"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
